$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Materialize the new block A5:C38 (and its style) in one shot: an
# explicit black font color on the whole block mints the new font
# (fontId 5) and cellXfs entry (index 7), and fills every row with
# styled A/B/C cells -- matching the row/column layout seen below.
$ws.Range("A5:C38").Font.Color = 0

# Re-home the old Phase labels and fill in the new task breakdown.
# Existing shared strings (Product Definition/Architecture/Realization/
# Bring Up/Test) are reused in place; everything else is new text that
# gets appended to the shared string table in this same order.
$ws.Range("A5").Value = "Product Definition"
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "Brainstorming"
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = "Research"
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = "Design Analysis"
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = "Flow Diagram "
$ws.Range("B10").Value = "Learning*"
$ws.Range("B11").Value = "Characteristics List"
$ws.Range("B12").Value = "Behavior Definition"
$ws.Range("B13").Value = "Interface List"
$ws.Range("A14").Value = "Architecture"
$ws.Range("B15").Value = "Source"
$ws.Range("B16").Value = "Major Compnents BOM"
$ws.Range("B17").Value = "HW Block Diagram"
$ws.Range("B18").Value = "Datasheet Research"
$ws.Range("B19").Value = "Product Architecture"
$ws.Range("C20").Value = "Calculations"
$ws.Range("C21").Value = "Current Budget "
$ws.Range("C22").Value = "Major Signals List"
$ws.Range("C23").Value = "Interface Pinouts"
$ws.Range("A24").Value = "Realization"
$ws.Range("B25").Value = "Prebuilt Eagle library Research"
$ws.Range("B26").Value = "Eagle library build"
$ws.Range("C27").Value = "Device"
$ws.Range("C28").Value = "Package "
$ws.Range("C29").Value = "Symbol"
$ws.Range("B30").Value = "Eagle schematic"
$ws.Range("C31").Value = "Place"
$ws.Range("C32").Value = "Wire"
$ws.Range("C33").Value = "Naming"
$ws.Range("C34").Value = "Error Checking"
$ws.Range("B35").Value = "Eagle Layout"
$ws.Range("A37").Value = "Bring Up"
$ws.Range("A38").Value = "Test "

# Match the recorded selection after the edit.
$ws.Range("A5:C38").Select()
